# Insert a new data row at row 41 (before the current "Madrigal / Primera / 44505" row),
# which pushes all subsequent rows (old 41..89) down by one (new 42..90).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's data.
$ws.Range("A41").Value2 = 5
$ws.Range("B41").Value2 = "Macroferia Regional de Talca"
$ws.Range("C41").Value2 = "Maule"
$ws.Range("D41").Value2 = 44781
$ws.Range("E41").Value2 = 7
$ws.Range("F41").Value2 = 100112013
$ws.Range("G41").Value2 = "Alcachofa"
$ws.Range("H41").Value2 = "Madrigal"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 300
$ws.Range("K41").Value2 = 14000
$ws.Range("L41").Value2 = 14000
$ws.Range("M41").Value2 = 14000
$ws.Range("N41").Value2 = "$/caja 40 unidades"
$ws.Range("O41").Value2 = "Provincia del Elquí"
$ws.Range("P41").Value2 = 350
$ws.Range("Q41").Value2 = 40
$ws.Range("R41").Value2 = "Hortaliza"
